$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,5).Value = 2
$ws.Cells.Item(2,6).Value = 0.6666666666666666
$ws.Cells.Item(2,7).Value = 0.1881506666666667
$ws.Cells.Item(2,8).Value = 0.564452
$ws.Cells.Item(2,9).Value = 0.06062261335217859
$ws.Cells.Item(2,10).Value = 0.0606226133521786
$ws.Cells.Item(2,13).Value = 8.554479333333333
$ws.Cells.Item(2,14).Value = 25.663438
$ws.Cells.Item(2,15).Value = 0.1655051910559175
$ws.Cells.Item(2,16).Value = 0.1655051910559175
$ws.Cells.Item(2,17).Value = 1.609530989552889
$ws.Cells.Item(2,18).Value = 14.485778905976
$ws.Cells.Item(2,19).Value = 0.01003335720516133
$ws.Cells.Item(2,20).Value = 0.01003335720516134

# Row 3
$ws.Cells.Item(3,5).Value = 2
$ws.Cells.Item(3,6).Value = 0.6666666666666666
$ws.Cells.Item(3,7).Value = 0.1881506666666667
$ws.Cells.Item(3,8).Value = 0.564452
$ws.Cells.Item(3,9).Value = 0.06062261335217859
$ws.Cells.Item(3,10).Value = 0.0606226133521786
$ws.Cells.Item(3,13).Value = 20.28486166666667
$ws.Cells.Item(3,14).Value = 60.854585
$ws.Cells.Item(3,15).Value = 0.392455200938143
$ws.Cells.Item(3,16).Value = 0.392455200938143
$ws.Cells.Item(3,17).Value = 3.816610245824445
$ws.Cells.Item(3,18).Value = 34.34949221242
$ws.Cells.Item(3,19).Value = 0.0237916599045246
$ws.Cells.Item(3,20).Value = 0.0237916599045246

# Row 4
$ws.Cells.Item(4,5).Value = 2
$ws.Cells.Item(4,6).Value = 0.6666666666666666
$ws.Cells.Item(4,7).Value = 0.1881506666666667
$ws.Cells.Item(4,8).Value = 0.564452
$ws.Cells.Item(4,9).Value = 0.06062261335217859
$ws.Cells.Item(4,10).Value = 0.0606226133521786
$ws.Cells.Item(4,13).Value = 5.037112666666666
$ws.Cells.Item(4,14).Value = 15.111338
$ws.Cells.Item(4,15).Value = 0.09745400763531942
$ws.Cells.Item(4,16).Value = 0.09745400763531943
$ws.Cells.Item(4,17).Value = 0.9477361063084444
$ws.Cells.Item(4,18).Value = 8.529624956775999
$ws.Cells.Item(4,19).Value = 0.005907916624496229
$ws.Cells.Item(4,20).Value = 0.005907916624496231

# Row 5
$ws.Cells.Item(5,5).Value = 2
$ws.Cells.Item(5,6).Value = 0.6666666666666666
$ws.Cells.Item(5,7).Value = 0.1881506666666667
$ws.Cells.Item(5,8).Value = 0.564452
$ws.Cells.Item(5,9).Value = 0.06062261335217859
$ws.Cells.Item(5,10).Value = 0.0606226133521786
$ws.Cells.Item(5,13).Value = 17.810622
$ws.Cells.Item(5,14).Value = 53.431866
$ws.Cells.Item(5,15).Value = 0.34458560037062
$ws.Cells.Item(5,16).Value = 0.34458560037062
$ws.Cells.Item(5,17).Value = 3.351080403047999
$ws.Cells.Item(5,18).Value = 30.159723627432
$ws.Cells.Item(5,19).Value = 0.02088967961799642
$ws.Cells.Item(5,20).Value = 0.02088967961799643

# Row 6
$ws.Cells.Item(6,7).Value = 2.915487666666667
$ws.Cells.Item(6,8).Value = 8.746463
$ws.Cells.Item(6,9).Value = 0.9393773866478214
$ws.Cells.Item(6,10).Value = 0.9393773866478214
$ws.Cells.Item(6,13).Value = 8.554479333333333
$ws.Cells.Item(6,14).Value = 25.663438
$ws.Cells.Item(6,15).Value = 0.1655051910559175
$ws.Cells.Item(6,16).Value = 0.1655051910559175
$ws.Cells.Item(6,17).Value = 24.94047899108822
$ws.Cells.Item(6,18).Value = 224.464310919794
$ws.Cells.Item(6,19).Value = 0.1554718338507562
$ws.Cells.Item(6,20).Value = 0.1554718338507562

# Row 7
$ws.Cells.Item(7,7).Value = 2.915487666666667
$ws.Cells.Item(7,8).Value = 8.746463
$ws.Cells.Item(7,9).Value = 0.9393773866478214
$ws.Cells.Item(7,10).Value = 0.9393773866478214
$ws.Cells.Item(7,13).Value = 20.28486166666667
$ws.Cells.Item(7,14).Value = 60.854585
$ws.Cells.Item(7,15).Value = 0.392455200938143
$ws.Cells.Item(7,16).Value = 0.392455200938143
$ws.Cells.Item(7,17).Value = 59.14026400920612
$ws.Cells.Item(7,18).Value = 532.262376082855
$ws.Cells.Item(7,19).Value = 0.3686635410336184
$ws.Cells.Item(7,20).Value = 0.3686635410336184

# Row 8
$ws.Cells.Item(8,7).Value = 2.915487666666667
$ws.Cells.Item(8,8).Value = 8.746463
$ws.Cells.Item(8,9).Value = 0.9393773866478214
$ws.Cells.Item(8,10).Value = 0.9393773866478214
$ws.Cells.Item(8,13).Value = 5.037112666666666
$ws.Cells.Item(8,14).Value = 15.111338
$ws.Cells.Item(8,15).Value = 0.09745400763531942
$ws.Cells.Item(8,16).Value = 0.09745400763531943
$ws.Cells.Item(8,17).Value = 14.68563985527711
$ws.Cells.Item(8,18).Value = 132.170758697494
$ws.Cells.Item(8,19).Value = 0.09154609101082319
$ws.Cells.Item(8,20).Value = 0.0915460910108232

# Row 9
$ws.Cells.Item(9,7).Value = 2.915487666666667
$ws.Cells.Item(9,8).Value = 8.746463
$ws.Cells.Item(9,9).Value = 0.9393773866478214
$ws.Cells.Item(9,10).Value = 0.9393773866478214
$ws.Cells.Item(9,13).Value = 17.810622
$ws.Cells.Item(9,14).Value = 53.431866
$ws.Cells.Item(9,15).Value = 0.34458560037062
$ws.Cells.Item(9,16).Value = 0.34458560037062
$ws.Cells.Item(9,17).Value = 51.926648776662
$ws.Cells.Item(9,18).Value = 467.339838989958
$ws.Cells.Item(9,19).Value = 0.3236959207526235
$ws.Cells.Item(9,20).Value = 0.3236959207526236
